# Insert a new price-report row (weekly update) above the current row 36,
# shifting the existing rows 36-48 down to 37-49, then populate the new
# row 36 with the latest "Haba" (Vega Modelo de Temuco) record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 36..48 down to 37..49 (carries formatting, incl. the date
# number format on column D, along with them).
$ws.Rows.Item(36).Insert()

# Fill in the new row 36 with the new weekly record.
$ws.Cells.Item(36, 1).Value = 10
$ws.Cells.Item(36, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(36, 3).Value = "La Araucanía"
$ws.Cells.Item(36, 4).Value = 44511
$ws.Cells.Item(36, 5).Value = 9
$ws.Cells.Item(36, 6).Value = 100112026
$ws.Cells.Item(36, 7).Value = "Haba"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 900
$ws.Cells.Item(36, 11).Value = 7000
$ws.Cells.Item(36, 12).Value = 8000
$ws.Cells.Item(36, 13).Value = 7556
$ws.Cells.Item(36, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Región del Maule"
$ws.Cells.Item(36, 16).Value = 302
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
